$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '26.914.11'
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").Value = '1.846.23'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.42'
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4708'
$ws.Range("E7").Value = '  +3.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3663'
$ws.Range("E8").Value = '  +1.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07154'
$ws.Range("E9").Value = '  +0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9276'
$ws.Range("E10").Value = '  +3.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.57'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07707'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").Value = '1.880.72'
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.283'
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.406'
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.33'
$ws.Range("E16").Value = '  +3.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008631'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").Value = '26.947.36'
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.43'
$ws.Range("E21").Value = '  +2.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.025'
$ws.Range("E22").Value = '  +1.31%  '
$ws.Range("E23").Value = '  +1.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.930'
$ws.Range("E24").Value = '  -1.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.95'
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.25'
$ws.Range("E26").Value = '  +2.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.015'
$ws.Range("E27").Value = '  -1.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.40'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.879'
$ws.Range("E29").Value = '  +0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08854'
$ws.Range("E30").Value = '  +1.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.213'
$ws.Range("E31").Value = '  +2.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.179'
$ws.Range("E32").Value = '  +6.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7487'
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.791'
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.477'
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.087'
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01940'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.967'
$ws.Range("E38").Value = '  +1.77%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05204'
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5228'
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.971'
$ws.Range("E41").Value = '  +2.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1517'
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.151'
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.49'
$ws.Range("E44").Value = '  +5.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4707'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.006'
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.79'
$ws.Range("E47").Value = '  +1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.597'
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '65.32'
$ws.Range("E49").Value = '  +2.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06036'
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8960'
$ws.Range("E51").Value = '  +5.69%  '
